# Files can be sent to scripts as arguments from de terminal
# Update the random bias/weight values on both layer sheets.

$wb = $excel.ActiveWorkbook

$ws0 = $wb.Worksheets.Item("Layer0")
$ws0.Range("B2").Value = 0.3494694828779881
$ws0.Range("C2").Value = -0.7815583031016935
$ws0.Range("B3").Value = -0.7728078242598293
$ws0.Range("C3").Value = 0.5587969675177827
$ws0.Range("B4").Value = -1.123941849698018
$ws0.Range("C4").Value = 0.04183512305456807

$ws1 = $wb.Worksheets.Item("Layer1")
$ws1.Range("B2").Value = -0.6422312916215966
$ws1.Range("C2").Value = -0.1278378689239725
$ws1.Range("B3").Value = -2.004016977517532
$ws1.Range("C3").Value = 0.9008433375932008
$ws1.Range("B4").Value = 0.6932072598580434
$ws1.Range("C4").Value = -0.5551909189675438
